# Auto-generated edit script: update crypto price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.384.78'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '3.612.09'
$ws.Range('E3').Value = '  +1.77%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '604.49'
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('D6').Value = '196.04'
$ws.Range('E6').Value = '  -0.76%  '
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -1.27%  '
$ws.Range('E10').Value = '  -1.03%  '
$ws.Range('D11').Value = '53.61'
$ws.Range('E11').Value = '  -1.00%  '
$ws.Range('E12').Value = '  -0.12%  '
$ws.Range('D13').Value = '9.57'
$ws.Range('E13').Value = '  -0.02%  '
$ws.Range('D14').Value = '4.181.90'
$ws.Range('E14').Value = '  +1.61%  '
$ws.Range('D15').Value = '13.00'
$ws.Range('E15').Value = '  +2.22%  '
$ws.Range('D16').Value = '596.28'
$ws.Range('E16').Value = '  -1.19%  '
$ws.Range('D17').Value = '70.451.44'
$ws.Range('E17').Value = '  +0.33%  '
$ws.Range('D18').Value = '3.613.32'
$ws.Range('E18').Value = '  +1.79%  '
$ws.Range('D19').Value = '19.04'
$ws.Range('E19').Value = '  -1.03%  '
$ws.Range('E20').Value = '  +1.43%  '
$ws.Range('E21').Value = '  -0.20%  '
$ws.Range('D22').Value = '17.88'
$ws.Range('E22').Value = '  -1.17%  '
$ws.Range('D23').Value = '5.19'
$ws.Range('E23').Value = '  -2.38%  '
$ws.Range('D24').Value = '101.81'
$ws.Range('E24').Value = '  -1.00%  '
$ws.Range('D25').Value = '4.64'
$ws.Range('E25').Value = '  +0.20%  '
$ws.Range('D26').Value = '3.03'
$ws.Range('E26').Value = '  -3.43%  '
$ws.Range('E27').Value = '  -1.70%  '
$ws.Range('D28').Value = '9.61'
$ws.Range('E28').Value = '  -0.55%  '
$ws.Range('D29').Value = '33.79'
$ws.Range('E29').Value = '  -0.27%  '
$ws.Range('D30').Value = '4.74'
$ws.Range('E30').Value = '  +6.36%  '
$ws.Range('E31').Value = '  +0.59%  '
$ws.Range('D32').Value = '12.30'
$ws.Range('E32').Value = '  -3.07%  '
$ws.Range('E33').Value = '  +1.35%  '
$ws.Range('D34').Value = '63.46'
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('D35').Value = '0.0₃0899'
$ws.Range('E35').Value = '  +6.75%  '
$ws.Range('D36').Value = '3.900.77'
$ws.Range('E36').Value = '  +2.92%  '
$ws.Range('D37').Value = '542.17'
$ws.Range('E37').Value = '  +11.12%  '
$ws.Range('E38').Value = '  +1.49%  '
$ws.Range('E39').Value = '  +0.11%  '
$ws.Range('D40').Value = '36.94'
$ws.Range('E40').Value = '  +0.34%  '
$ws.Range('E41').Value = '  -1.42%  '
$ws.Range('E42').Value = '  -4.53%  '
$ws.Range('D43').Value = '0.134'
$ws.Range('E43').Value = '  -1.19%  '
$ws.Range('E44').Value = '  -0.60%  '
$ws.Range('E45').Value = '  +4.50%  '
$ws.Range('D46').Value = '2.87'
$ws.Range('E46').Value = '  +0.46%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').Value = '8.61'
$ws.Range('E48').Value = '  -1.04%  '
$ws.Range('E49').Value = '  -0.20%  '
$ws.Range('D50').Value = '0.000252'
$ws.Range('E50').Value = '  +0.65%  '
$ws.Range('D51').Value = '1.31'
$ws.Range('E51').Value = '  -0.24%  '
